$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new header cells F1 (Center_X) and G1 (Center_Z), matching the
# --- style of the existing header cells (e.g. E1) ---
$ws.Range("E1").Copy()
$ws.Range("F1:G1").PasteSpecial(-4122)
$ws.Range("F1").Value = "Center_X"
$ws.Range("G1").Value = "Center_Z"

# --- Extend the A-column "id" style (bold/border) down to the two new rows ---
$ws.Range("A5").Copy()
$ws.Range("A6:A7").PasteSpecial(-4122)

# --- Fill in the id column for the new rows ---
$ws.Range("A6").Value = 50
$ws.Range("A7").Value = 60

# --- Update data rows 2-7 for columns B..G with the new values ---
$weights = 85979.89872084881
$momentX = 2629653.239222029
$momentZ = 57615.16669750272
$momentY = 0
$centerX = 30.58451194225911
$centerZ = 0.6701004252698883

foreach ($r in 2..7) {
    $ws.Cells.Item($r, 2).Value = $weights
    $ws.Cells.Item($r, 3).Value = $momentX
    $ws.Cells.Item($r, 4).Value = $momentZ
    $ws.Cells.Item($r, 5).Value = $momentY
    $ws.Cells.Item($r, 6).Value = $centerX
    $ws.Cells.Item($r, 7).Value = $centerZ
}

Write-Host "Gravity centers added"
